$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: classical-best-embeddings vs. classical-best-tfidf -> classical-best-embed vs. classical-best-tfidf
$ws.Range("A2").Value = "classical-best-embed vs. classical-best-tfidf"
$ws.Range("C2").Value = 0.057
$ws.Range("E2").Value = 0.028
$ws.Range("F2").Value = 0.023
$ws.Range("H2").Value = 0.041
$ws.Range("I2").Value = 0.032
$ws.Range("J2").Value = 0.04

# Row 3: BERT-base vs. classical-best-tfidf
$ws.Range("C3").Value = 0.064
$ws.Range("D3").Value = 0.08
$ws.Range("E3").Value = 0.096
$ws.Range("F3").Value = 0.073
$ws.Range("G3").Value = 0.13
$ws.Range("H3").Value = 0.106
$ws.Range("I3").Value = 0.078

# Row 4: BERT-base vs. classical-best-embeddings -> BERT-base vs. classical-best-embed
$ws.Range("A4").Value = "BERT-base vs. classical-best-embed"
$ws.Range("C4").Value = 0.007
$ws.Range("D4").Value = 0.059
$ws.Range("E4").Value = 0.068
$ws.Range("F4").Value = 0.05
$ws.Range("G4").Value = 0.06
$ws.Range("H4").Value = 0.065
$ws.Range("I4").Value = 0.046
$ws.Range("J4").Value = 0.052

# Row 5: BERT-base-nli vs. classical-best-tfidf
$ws.Range("B5").Value = 0.46
$ws.Range("C5").Value = 0.156
$ws.Range("D5").Value = 0.087
$ws.Range("E5").Value = 0.094
$ws.Range("G5").Value = 0.073
$ws.Range("H5").Value = 0.049
$ws.Range("I5").Value = 0.098
$ws.Range("J5").Value = 0.086

# Row 6: BERT-base-nli vs. classical-best-embeddings -> BERT-base-nli vs. classical-best-embed
$ws.Range("A6").Value = "BERT-base-nli vs. classical-best-embed"
$ws.Range("B6").Value = 0.46
$ws.Range("C6").Value = 0.099
$ws.Range("D6").Value = 0.066
$ws.Range("E6").Value = 0.066
$ws.Range("F6").Value = 0.032
$ws.Range("G6").Value = 0.003
$ws.Range("H6").Value = 0.008
$ws.Range("I6").Value = 0.066
$ws.Range("J6").Value = 0.046

# Row 7: BERT-base-nli vs. BERT-base
$ws.Range("B7").Value = 0.46
$ws.Range("C7").Value = 0.092
$ws.Range("D7").Value = 0.007
$ws.Range("E7").Value = -0.002
$ws.Range("F7").Value = -0.018
$ws.Range("G7").Value = -0.057
$ws.Range("H7").Value = -0.057
$ws.Range("I7").Value = 0.02
$ws.Range("J7").Value = -0.006
